$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (LinearRegression) - update values
$ws.Range("B2").Value = 5489827817512978
$ws.Range("C2").Value = 5489827817512978
$ws.Range("D2").Value = 5489827817512978

# Row 3 (RandomForestRegressor) - update values
$ws.Range("B3").Value = 9942265647837.883
$ws.Range("C3").Value = 7944078720463.133
$ws.Range("D3").Value = 7016406959225.904

# Row 4 (GradientBoostingRegressor -> DecisionTreeRegressor)
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 4958455199013.18
$ws.Range("C4").Value = 4958455199013.18
$ws.Range("D4").Value = 4957675996498.995

# Row 5 (AdaBoostRegressor -> MLPRegressor)
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 151260026024762.9
$ws.Range("C5").Value = 199577340389577.3
$ws.Range("D5").Value = 253407287789214.3
